# Auto-generated edit script: updates profit-calculation columns (H-N)
# on rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching a
# scheduled-runner refresh of market-price-derived figures.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2670.5
$ws.Range("I74").Value = 2239.6667
$ws.Range("J74").Value = 2929
$ws.Range("K74").Value = 2239.6667
$ws.Range("L74").Value = 2929
$ws.Range("M74").Value = -1303.6667
$ws.Range("N74").Value = -4801
$ws.Range("H76").Value = 11912912
$ws.Range("I76").Value = 8781.883
$ws.Range("J76").Value = 30310204
$ws.Range("K76").Value = 8781.883
$ws.Range("L76").Value = 30310204
$ws.Range("M76").Value = -8466.883
$ws.Range("N76").Value = -30310834
$ws.Range("H77").Value = 2670.5
$ws.Range("I77").Value = 2239.6667
$ws.Range("J77").Value = 2929
$ws.Range("K77").Value = 11198.3335
$ws.Range("L77").Value = 14645
$ws.Range("M77").Value = -6518.333500000001
$ws.Range("N77").Value = -24005
$ws.Range("H79").Value = 11912912
$ws.Range("I79").Value = 8781.883
$ws.Range("J79").Value = 30310204
$ws.Range("K79").Value = 8781.883
$ws.Range("L79").Value = 30310204
$ws.Range("M79").Value = -7689.883
$ws.Range("N79").Value = -30312388
$ws.Range("H111").Value = 973.55554
$ws.Range("I111").Value = 482.5
$ws.Range("J111").Value = 1366.4
$ws.Range("K111").Value = 1447.5
$ws.Range("L111").Value = 4099.200000000001
$ws.Range("M111").Value = 1619.5
$ws.Range("N111").Value = -10233.2
$ws.Range("H113").Value = 3184.2083
$ws.Range("I113").Value = 5201.25
$ws.Range("J113").Value = 2780.8
$ws.Range("K113").Value = 5201.25
$ws.Range("L113").Value = 2780.8
$ws.Range("M113").Value = -1947.25
$ws.Range("N113").Value = -9288.799999999999
$ws.Range("H131").Value = 1488.3684
$ws.Range("I131").Value = 441.8125
$ws.Range("J131").Value = 7070
$ws.Range("K131").Value = 1325.4375
$ws.Range("L131").Value = 21210
$ws.Range("M131").Value = 3714.5625
$ws.Range("N131").Value = -31290
$ws.Range("H137").Value = 1995.2543
$ws.Range("I137").Value = 1960.85
$ws.Range("J137").Value = 2067.6843
$ws.Range("K137").Value = 5882.549999999999
$ws.Range("L137").Value = 6203.0529
$ws.Range("M137").Value = -3332.549999999999
$ws.Range("N137").Value = -11303.0529
$ws.Range("H138").Value = 3917.4036
$ws.Range("I138").Value = 1387.3
$ws.Range("J138").Value = 9870.588
$ws.Range("K138").Value = 4161.9
$ws.Range("L138").Value = 29611.764
$ws.Range("M138").Value = 978.1000000000004
$ws.Range("N138").Value = -39891.764

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3760.3333
$ws.Range("I110").Value = 3488.875
$ws.Range("J110").Value = 4303.25
$ws.Range("K110").Value = 3488.875
$ws.Range("L110").Value = 4303.25
$ws.Range("M110").Value = -1443.875
$ws.Range("N110").Value = -8393.25

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 901311.25
$ws.Range("I22").Value = 1039836.06
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 1039836.06
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -1039663.06
$ws.Range("N22").Value = -1246
$ws.Range("H107").Value = 35716720
$ws.Range("I107").Value = 100001380
$ws.Range("J107").Value = 3025.111
$ws.Range("K107").Value = 100001380
$ws.Range("L107").Value = 3025.111
$ws.Range("M107").Value = -99999460
$ws.Range("N107").Value = -6865.111
$ws.Range("H134").Value = 3478319.5
$ws.Range("I134").Value = 2968.2173
$ws.Range("J134").Value = 12359773
$ws.Range("K134").Value = 8904.651899999999
$ws.Range("L134").Value = 37079319
$ws.Range("M134").Value = -6369.651899999999
$ws.Range("N134").Value = -37084389

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1104.25
$ws.Range("I31").Value = 893.6842
$ws.Range("J31").Value = 1225.4849
$ws.Range("K31").Value = 893.6842
$ws.Range("L31").Value = 1225.4849
$ws.Range("M31").Value = -598.6842
$ws.Range("N31").Value = -1815.4849
$ws.Range("H34").Value = 1104.25
$ws.Range("I34").Value = 893.6842
$ws.Range("J34").Value = 1225.4849
$ws.Range("K34").Value = 893.6842
$ws.Range("L34").Value = 1225.4849
$ws.Range("M34").Value = -691.6842
$ws.Range("N34").Value = -1629.4849
$ws.Range("H107").Value = 494.7143
$ws.Range("I107").Value = 489.72726
$ws.Range("J107").Value = 513
$ws.Range("K107").Value = 489.72726
$ws.Range("L107").Value = 513
$ws.Range("M107").Value = 1430.27274
$ws.Range("N107").Value = -4353

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 3013
$ws.Range("I101").Value = 3026
$ws.Range("J101").Value = 3000
$ws.Range("K101").Value = 9078
$ws.Range("L101").Value = 9000
$ws.Range("M101").Value = -6644
$ws.Range("N101").Value = -13868
$ws.Range("H109").Value = 5999
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H131").Value = 810.8200000000001
$ws.Range("I131").Value = 504.45456
$ws.Range("J131").Value = 848.6853599999999
$ws.Range("K131").Value = 1513.36368
$ws.Range("L131").Value = 2546.05608
$ws.Range("M131").Value = 3526.63632
$ws.Range("N131").Value = -12626.05608

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 1450.5
$ws.Range("I17").Value = 401
$ws.Range("J17").Value = 2500
$ws.Range("K17").Value = 401
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = -233
$ws.Range("N17").Value = -2836
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H31").Value = 1500
$ws.Range("I31").Value = 1500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1208
$ws.Range("H37").Value = 1500
$ws.Range("I37").Value = 1500
$ws.Range("K37").Value = 1500
$ws.Range("M37").Value = -1223
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H107").Value = 4684.913
$ws.Range("I107").Value = 291.375
$ws.Range("J107").Value = 7028.1333
$ws.Range("K107").Value = 291.375
$ws.Range("L107").Value = 7028.1333
$ws.Range("M107").Value = 1628.625
$ws.Range("N107").Value = -10868.1333
$ws.Range("H132").Value = 4669.84
$ws.Range("I132").Value = 2411.3142
$ws.Range("J132").Value = 9939.733
$ws.Range("K132").Value = 7233.942599999999
$ws.Range("L132").Value = 29819.199
$ws.Range("M132").Value = -4703.942599999999
$ws.Range("N132").Value = -34879.199

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 14708471
$ws.Range("I40").Value = 2083.9167
$ws.Range("K40").Value = 2083.9167
$ws.Range("M40").Value = -1947.9167
$ws.Range("H122").Value = 5341.769
$ws.Range("I122").Value = 6354.05
$ws.Range("K122").Value = 19062.15
$ws.Range("M122").Value = -16612.15
$ws.Range("H136").Value = 34878264
$ws.Range("I136").Value = 50905924
$ws.Range("J136").Value = 23811544
$ws.Range("K136").Value = 152717772
$ws.Range("L136").Value = 71434632
$ws.Range("M136").Value = -152715222
$ws.Range("N136").Value = -71439732

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H122").Value = 35859.668
$ws.Range("I122").Value = 57166.89
$ws.Range("K122").Value = 171500.67
$ws.Range("M122").Value = -169050.67
$ws.Range("H126").Value = 2553.6667
$ws.Range("I126").Value = 2297.2144
$ws.Range("J126").Value = 3451.25
$ws.Range("K126").Value = 6891.6432
$ws.Range("L126").Value = 10353.75
$ws.Range("M126").Value = -4421.6432
$ws.Range("N126").Value = -15293.75
$ws.Range("H132").Value = 29883.691
$ws.Range("I132").Value = 34919.656
$ws.Range("J132").Value = 13768.6
$ws.Range("K132").Value = 104758.968
$ws.Range("L132").Value = 41305.8
$ws.Range("M132").Value = -102228.968
$ws.Range("N132").Value = -46365.8
$ws.Range("H136").Value = 8930097
$ws.Range("I136").Value = 14286499
$ws.Range("K136").Value = 42859497
$ws.Range("M136").Value = -42856947
